$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roles")

# Insert two new rows above the old "header repeat" row (old row 34 -> becomes row 36)
$ws.Rows("34:35").Insert()

# Row 34 only uses columns A, B, E - clear out the inherited column styling on C/D
$ws.Range("C34:D34").Clear()

$ws.Range("A34").Value = "Cursed"
$ws.Range("B34").Value = "Village (until dead, then werewolf)"
$ws.Range("E34").Value = -3

# Row 35 uses columns A-E
$ws.Range("A35").Value = "Old Man"
$ws.Range("B35").Value = "Village"
$ws.Range("C35").Value = "Dies on night (# of Wolves) + 1"
$ws.Range("D35").Value = "Never"
$ws.Range("E35").Value = 1

# Give the new rows the same row height behavior as the rest of the sheet
$ws.Rows("34:35").RowHeight = 15

# Column B width change (target stored width 17.140625; engine quantizes
# ColumnWidth to its internal pixel grid, so aim for the nearest achievable value)
$ws.Columns("B").ColumnWidth = 16.25

# Update the active selection
$ws.Range("E35").Select()

Write-Output "done"
